$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample Types")
$ws.Activate() | Out-Null

# --- Workbook-level metadata ---
# (absPath / yWindow are Excel UI/machine-specific metadata not exposed via
# the Excel object model; they cannot be set through COM automation.)

# --- New column E: "SampleType Processing Notes" ---
$ws.Range("E1").Value = "SampleType Processing Notes"
$ws.Range("E3").Value = "use Not Recorded"
$ws.Range("E8").Value = "update MatrixCode <-runoff, if SAMP_TYPE_CD =9"
$ws.Range("E9").Value = "update MatrixCode <-runoff, if SAMP_TYPE_CD =10"

# Match the style already used on the other descriptive columns (C/D -> style index 2,
# font "Segoe UI" size 11)
foreach ($addr in @("E1", "E3", "E8", "E9")) {
  $ws.Range($addr).Font.Name = "Segoe UI"
  $ws.Range($addr).Font.Size = 11
}

# Set column E width as close as this runtime allows to the authored 47.7109375
$ws.Columns.Item(5).ColumnWidth = 46.8

# --- Selection / view state ---
$ws.Range("D17").Select() | Out-Null
